# New crime data collected - weekly CompStat update (1/30/2023 - 2/5/2023)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates: report "Volume 30 Number 4" -> "Number 5" and the
# week-covering date range 1/23/2023-1/29/2023 -> 1/30/2023-2/5/2023.
# Use Characters() so only the affected substring of the rich-text cell
# is touched (other runs / formatting are left alone).
# ---------------------------------------------------------------------
$ws.Range("A8").Characters(21, 1).Text = "5"
$ws.Range("C9").Characters(27, 9).Text = "1/30/2023"
$ws.Range("C9").Characters(47, 9).Text = "2/5/2023"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("L15").Value = -50
$ws.Range("N15").Value = -66.666666666666

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = 48.148148148148
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 37
$ws.Range("K16").Value = 59.459459459459
$ws.Range("L16").Value = 103.448275862069
$ws.Range("M16").Value = 227.777777777778
$ws.Range("N16").Value = -78.066914498141

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 123.529411764706
$ws.Range("I17").Value = 56
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = 115.384615384615
$ws.Range("L17").Value = 60
$ws.Range("M17").Value = 124
$ws.Range("N17").Value = -17.647058823529

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 35
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = -7.894736842105
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = 10.416666666666
$ws.Range("L18").Value = 140.909090909091
$ws.Range("M18").Value = 15.217391304347
$ws.Range("N18").Value = -83.012820512820

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 36
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 161
$ws.Range("G19").Value = 119
$ws.Range("H19").Value = 35.294117647058
$ws.Range("I19").Value = 237
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = 39.411764705882
$ws.Range("L19").Value = 146.875
$ws.Range("M19").Value = -2.066115702479
$ws.Range("N19").Value = -78.116343490304

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 20
$ws.Range("L20").Value = -14.285714285714
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -87.234042553191

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 5.357142857142
$ws.Range("F21").Value = 281
$ws.Range("G21").Value = 206
$ws.Range("H21").Value = 36.407766990291
$ws.Range("I21").Value = 412
$ws.Range("J21").Value = 287
$ws.Range("K21").Value = 43.554006968641
$ws.Range("L21").Value = 114.583333333333
$ws.Range("M21").Value = 23.353293413173
$ws.Range("N21").Value = -76.905829596412

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 3
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = 36.363636363636
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 15
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -14.285714285714

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 72
$ws.Range("D24").Value = 47
$ws.Range("E24").Value = 53.191489361702
$ws.Range("F24").Value = 276
$ws.Range("G24").Value = 168
$ws.Range("H24").Value = 64.285714285714
$ws.Range("I24").Value = 360
$ws.Range("J24").Value = 217
$ws.Range("K24").Value = 65.898617511520
$ws.Range("L24").Value = 94.594594594594
$ws.Range("M24").Value = -21.397379912663

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("E25").Value = -37.5
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 54
$ws.Range("H25").Value = 11.111111111111
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 71
$ws.Range("K25").Value = 16.901408450704
$ws.Range("L25").Value = 13.698630136986
$ws.Range("M25").Value = 53.703703703703

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("L26").Value = -50

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 300
$ws.Range("F27").Value = 15
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 21
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 200

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic. : L/N go from "***.*" placeholder text to real
# percentage values now that there is prior-period data to compare.
# (N22 is itself a "***.*" text placeholder, so borrow the number format
# from N16, a numeric percentage cell in the same N column.)
# ---------------------------------------------------------------------
$ws.Range("L28").Value = -100
$ws.Range("L28").NumberFormat = $ws.Range("L22").NumberFormat
$ws.Range("N28").Value = -100
$ws.Range("N28").NumberFormat = $ws.Range("N16").NumberFormat

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc. : same treatment as row 28.
# ---------------------------------------------------------------------
$ws.Range("L29").Value = -100
$ws.Range("L29").NumberFormat = $ws.Range("L22").NumberFormat
$ws.Range("N29").Value = -100
$ws.Range("N29").NumberFormat = $ws.Range("N16").NumberFormat

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes : C/F go from "0" placeholder text to real counts.
# ---------------------------------------------------------------------
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = $ws.Range("C16").NumberFormat
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = $ws.Range("F16").NumberFormat
$ws.Range("I30").Value = 2
